# Update database and change read_price algorithm
# - Drop oldest period (1396/12) and its publish date, shift all periods/dates
#   one column to the left (D<-E, E<-F, F<-G, G<-H), and add the newest
#   period (1401/12) with its publish date in column H.
# - Update all the financial figures accordingly (same left-shift + new
#   column H figures taken from the latest filing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---- Header row 8: financial period labels ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Header row 9: publish dates ----
$ws.Range("D9").Value = "1399-03-20 (8)"
$ws.Range("E9").Value = "1400-04-02 (8)"
$ws.Range("F9").Value = "1401-04-08 (8)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30"

# ---- Row 11: فروش (Sales) ----
$ws.Range("D11").Value = 730007
$ws.Range("E11").Value = 608284
$ws.Range("F11").Value = 633495
$ws.Range("G11").Value = 730059
$ws.Range("H11").Value = 607026

# ---- Row 12: بهای تمام شده کالای فروش رفته (COGS) ----
$ws.Range("D12").Value = -294919
$ws.Range("E12").Value = -275826
$ws.Range("F12").Value = -271312
$ws.Range("G12").Value = -554271
$ws.Range("H12").Value = -417553

# ---- Row 13: سود (زیان) ناخالص (Gross profit) ----
$ws.Range("D13").Value = 435088
$ws.Range("E13").Value = 332458
$ws.Range("F13").Value = 362183
$ws.Range("G13").Value = 175787
$ws.Range("H13").Value = 189472

# ---- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ----
$ws.Range("D14").Value = -103793
$ws.Range("E14").Value = -129214
$ws.Range("F14").Value = -150379
$ws.Range("G14").Value = -124379
$ws.Range("H14").Value = -131205

# Row 15 (impairment expense) stays all "-" (unchanged)

# ---- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other op. income/exp) ----
$ws.Range("D16").Value = 31531
$ws.Range("E16").Value = 29420
$ws.Range("F16").Value = 43214
$ws.Range("G16").Value = -2613
$ws.Range("H16").Value = 30718

# ---- Row 17: سود (زیان) عملیاتی (Operating profit) ----
$ws.Range("D17").Value = 362826
$ws.Range("E17").Value = 232664
$ws.Range("F17").Value = 255018
$ws.Range("G17").Value = 48795
$ws.Range("H17").Value = 88985

# ---- Row 18: هزینه های مالی (Financial expenses) -- now all "-" ----
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = "-"

# ---- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-op. income/exp) ----
$ws.Range("D19").Value = 7848
$ws.Range("E19").Value = -3514
$ws.Range("F19").Value = 1094
$ws.Range("G19").Value = 14164
$ws.Range("H19").Value = 625

# ---- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ----
$ws.Range("D20").Value = 370673
$ws.Range("E20").Value = 229150
$ws.Range("F20").Value = 256112
$ws.Range("G20").Value = 62959
$ws.Range("H20").Value = 89610

# ---- Row 21: مالیات (Tax) ----
$ws.Range("D21").Value = -1939
$ws.Range("E21").Value = -58
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = "-"
$ws.Range("H21").Value = "-"

# ---- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops) ----
$ws.Range("D22").Value = 368734
$ws.Range("E22").Value = 229092
$ws.Range("F22").Value = 256112
$ws.Range("G22").Value = 62959
$ws.Range("H22").Value = 89610

# Row 23 (discontinued ops) stays all "-" (unchanged)

# ---- Row 24: سود (زیان) خالص (Net profit) ----
$ws.Range("D24").Value = 368734
$ws.Range("E24").Value = 229092
$ws.Range("F24").Value = 256112
$ws.Range("G24").Value = 62959
$ws.Range("H24").Value = 89610

# Row 25 (EPS after tax) stays all 0 (unchanged)

# ---- Row 26: سرمایه (Capital) ----
$ws.Range("D26").Value = 23725
$ws.Range("E26").Value = 18709
$ws.Range("F26").Value = 10615
$ws.Range("G26").Value = 9096
$ws.Range("H26").Value = 6801

# Row 27 (EPS based on latest capital) stays all 0 (unchanged)
